$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999998820958025,
    0.99473536047671907,
    0.97267663198450971,
    0.96228009345247245,
    0.95227538980656767,
    0.92779046805842436,
    0.92639017159872794,
    0.92571233365305694,
    0.92824024997586485,
    0.93176091471006162,
    0.93247374157513396,
    0.93116045181563911,
    0.92613098938267657,
    0.92448299492901742,
    0.91776192070965457,
    0.9105760265671945,
    0.90686827019896166,
    0.90575939275534445,
    0.99368223474363382,
    0.98656545449798849,
    0.98516701411980656,
    0.98390251718543975,
    0.96564111267420383,
    0.95261948142355646,
    0.94616237828032035,
    0.92413363778223356,
    0.91928525255916194,
    0.89780640067994322,
    0.88253308995658297,
    0.8759623912966773,
    0.86830861112804991,
    0.86662935718447653,
    0.86610937350690631
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
